$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top - shifts all existing data down by one row.
$ws.Rows.Item(1).Insert()

# New header row: "Russian" / "English" labels above each of the three
# Russian/English column pairs (A:C, E:G, I:K), leaving the spacer
# column (B, F, J) blank.
$headerRange = $ws.Range("A1,C1,E1,G1,I1,K1")
$headerRange.Value = "Russian"

$ws.Range("A1").Value = "Russian"
$ws.Range("C1").Value = "English"
$ws.Range("E1").Value = "Russian"
$ws.Range("G1").Value = "English"
$ws.Range("I1").Value = "Russian"
$ws.Range("K1").Value = "English"

$headerRange.Font.Bold = $true
$headerRange.Font.Size = 11
$headerRange.Font.Color = 255

$ws.Range("K16").Select()
$ws.Range("I1:K1").Select()
